# LOM3219.docx restructuring edit
#
# The document's sections got reshuffled: several blocks of run text were
# cut from one paragraph and pasted into another, while every paragraph's
# own formatting container (style / bold / italic run-properties) stayed
# exactly where it was. So for every affected paragraph we only rewrite
# the *text payload* (using a vertical-tab char, [char]11, wherever a
# <w:br/> line break belongs — Word's Range.Text turns that into <w:br/>
# automatically), leaving paragraph count/order/styles untouched.

$d = $word.ActiveDocument
$VB = [char]11   # vertical tab -> becomes <w:br/> when assigned to Range.Text

# ---- literal text blocks reused below -------------------------------------------------

$PT_OBJETIVOS = "A maturação da nanotecnologia revelou que se trata de uma disciplina única e distinta, em vez de uma especialização dentro de um campo maior. Um curso sobre esse assunto envolve química, física e engenharia focada em Nano. Deve ser integrado, multidisciplinar e especificamente em Nano. A ideia é construir uma base sólida nos métodos de caracterização e fabricação enquanto integra a físicas e a química relevantes aos problemas envolvidos. Examinando os aspectos de engenharia, bem como nanomateriais e aplicações específicas nos setores de energia e eletrônica."

$EN_OBJETIVOS = "Nanotechnology maturation has revealed that it is a unique and distinct discipline rather than a specialization within a larger field. A course on this subject involves chemistry, physics and engineering focused on Nano. It must be integrated, multidisciplinary and specifically in nano. The idea is to build a solid foundation on characterization and manufacturing methods while integrating with physical and chemistry relevant to the problems involved. Examining engineering aspects as well as nanomaterials and specific applications in the energy and electronics sectors."

$PT_RESUMIDO = "Nanociência e nanotecnologia: princípios e aplicações."
$EN_RESUMIDO = "Nanoscience and nanotechnology: principles and applications."

$PT_PROGRAMA = "Perspectivas: nanociência e nanotecnologia - a distinção; Implicações sociais de nano" + $VB + `
    "Nanotools: métodos de caracterização; Métodos de fabricação" + $VB + `
    "Física: Propriedades e fenômenos: materiais, estrutura e nanosurface; Energia na nanoescala" + $VB + `
    "Química: síntese e modificação: nanomateriais à base de carbono; Interações químicas na nanoescala" + $VB + `
    "Aplicações: nanoetronics; nanomagnetismo; nanomecânica"

$EMERSON = "7290967 - Emerson Gonçalves de Melo"
$LUIZ = "1176388 - Luiz Tadeu Fernandes Eleno"

$AULAS_CRITERIO = "Aulas expositivas e seminários." + $VB + `
    "Critério" + $VB + `
    "Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3" + $VB + `
    "Norma de Recuperação" + $VB + `
    "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

$DUAS_PROVAS = "Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3"

$APLICACAO = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

$BIBLIO = "Gabor L. Hornyak, H.F. Tibbals, Joydeep Dutta, John J. Moore. Introduction to Nanoscience and Nanotechnology. CRC Press. 2009" + $VB + `
    "TIMP, G. Nanotechnology, Springer, 1998." + $VB + `
    "Bhushan, B. (ed.) Springer Handbook of Nanotechnology, Springer, 2010."

# ---- helper: replace the full text of a whole paragraph (by 1-based index) -----------
function Set-ParagraphText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $rng = $d.Range($p.Range.Start, $p.Range.End)
    $rng.Text = $newText
}

# Process paragraphs from the bottom of the document upward so that the
# Start/End offsets of paragraphs we haven't edited yet stay valid.

# Paragraph 19: Bibliografia body -> now just the "Luiz" docente line
Set-ParagraphText 19 $LUIZ

# Paragraph 17: Avaliação bullet list — three bold labels ("Método: ",
# "Critério: ", "Norma de recuperação: ") stay put; only the plain runs
# between/after them change. Each sub-edit resizes the paragraph, so the
# paragraph's Range (and every anchor inside it) is re-fetched fresh
# *after* every write instead of being reused from stale Int64 offsets.

# Segment after "Norma de recuperação: " runs to the end of the paragraph.
$p17 = $d.Paragraphs.Item(17)
$lblNorma = $d.Range($p17.Range.Start, $p17.Range.End)
$lblNorma.Find.Execute("Norma de recuperação: ") | Out-Null
$d.Range($lblNorma.End, $p17.Range.End).Text = $EMERSON

# Segment between "Critério: " and "Norma de recuperação: ".
$p17 = $d.Paragraphs.Item(17)
$lblCriterio = $d.Range($p17.Range.Start, $p17.Range.End)
$lblCriterio.Find.Execute("Critério: ") | Out-Null
$p17 = $d.Paragraphs.Item(17)
$lblNorma2 = $d.Range($p17.Range.Start, $p17.Range.End)
$lblNorma2.Find.Execute("Norma de recuperação: ") | Out-Null
$d.Range($lblCriterio.End, $lblNorma2.Start).Text = $BIBLIO + $VB

# Segment between "Método: " and "Critério: ".
$p17 = $d.Paragraphs.Item(17)
$lblMetodo = $d.Range($p17.Range.Start, $p17.Range.End)
$lblMetodo.Find.Execute("Método: ") | Out-Null
$p17 = $d.Paragraphs.Item(17)
$lblCriterio2 = $d.Range($p17.Range.Start, $p17.Range.End)
$lblCriterio2.Find.Execute("Critério: ") | Out-Null
$d.Range($lblMetodo.End, $lblCriterio2.Start).Text = $APLICACAO + $VB

# Paragraph 14: "Programa" detail (PT) -> shrinks down to just the exam-grade line
Set-ParagraphText 14 $DUAS_PROVAS

# Paragraph 12: EN objectives text (italic) moves here from paragraph 7
Set-ParagraphText 12 $EN_OBJETIVOS

# Paragraph 11: "Programa resumido" body -> the old "Método" run content
Set-ParagraphText 11 $AULAS_CRITERIO

# Paragraph 9: "Docente(s)" bullet list -> PT objectives text + PT Programa detail
$PT_DOCENTE_PARA = $PT_OBJETIVOS + $VB + $PT_PROGRAMA
Set-ParagraphText 9 $PT_DOCENTE_PARA

# Paragraph 7: EN objectives (italic) -> EN "Programa resumido" short line
Set-ParagraphText 7 $EN_RESUMIDO

# Paragraph 6: PT objectives text -> PT "Programa resumido" short line
Set-ParagraphText 6 $PT_RESUMIDO

Write-Host "Done."
